# Auto-generated edit script: update crypto Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a pure numeric-looking string but must remain
# stored as text (matching the source inlineStr cells) -- force text format
# before assignment, then restore the default style so no stray formatting
# is introduced.
$textForceCells = @("D4", "D5", "D6", "D7", "D12", "D13", "D15", "D17", "D20", "D21", "D23", "D24", "D25", "D27", "D28", "D30", "D31", "D32", "D37", "D38", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.156.00"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.526.53"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "537.29"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "137.35"
$ws.Range("E6").Value = "  -1.45%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("D9").Value = "2.523.40"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "0.349"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "2.971.78"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "23.08"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "59.077.43"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "2.522.94"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "4.30"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "324.03"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "5.96"
$ws.Range("E23").Value = "  +1.97%  "
$ws.Range("D24").Value = "65.97"
$ws.Range("E24").Value = "  +4.14%  "
$ws.Range("D25").Value = "0.424"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "7.56"
$ws.Range("E28").Value = "  -3.21%  "
$ws.Range("D29").Value = "0.0₃0777"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "6.72"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").Value = "1.78"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "166.80"
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("E33").Value = "  +5.68%  "
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "4.12"
$ws.Range("E37").Value = "  -2.29%  "
$ws.Range("D38").Value = "1.55"
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "0.817"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").Value = "3.63"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").Value = "285.27"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("D44").Value = "132.52"
$ws.Range("E44").Value = "  +6.07%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "0.606"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").Value = "10.90"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "0.0926"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D49").Value = "0.0510"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("E50").Value = "  -1.34%  "
$ws.Range("D51").Value = "17.41"
$ws.Range("E51").Value = "  -2.42%  "

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
